$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19.32001113891602
$ws.Range("C3").Value = 17.85588264465332
$ws.Range("C4").Value = 17.46892929077148
$ws.Range("C5").Value = 17.57693290710449
$ws.Range("C6").Value = 19.30093765258789
